$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2-3: account holder name and card number
$ws.Range("C2").Value = "Hartmut"
# Card number must stay text (16-digit number, too long to be a safe number).
# Use a text-producing formula, then collapse it down to a plain value so the
# cell keeps its original style and becomes a normal (non-formula) text cell.
$ws.Range("B3").Formula = "=""2570314725427075"""
$ws.Range("B3").Copy()
$ws.Range("B3").PasteSpecial(-4163) # xlPasteValues
$ws.Range("C3").Value = "Mohaupt"

# Opening balance line
$ws.Range("D5").Value = "KONTOSTAND AM 21.01.2024"

# Row 6
$ws.Range("B6").Value = "22.01."
$ws.Range("C6").Value = "23.01."
$ws.Range("D6").Value = "RECHNUNG VODAFONE GMBH 28076163"
$ws.Range("E6").Value = "39,11-"

# Row 7
$ws.Range("B7").Value = "25.01."
$ws.Range("C7").Value = "26.01."
$ws.Range("D7").Value = "PAYPAL VTNIJO"
$ws.Range("E7").Value = "74,25-"

# Row 8
$ws.Range("B8").Value = "27.01."
$ws.Range("C8").Value = "28.01."
$ws.Range("D8").Value = "KARTENZAHLUNG JET TANKSTELLE"
$ws.Range("E8").Value = "59,95-"

# Row 9
$ws.Range("B9").Value = "30.01."
$ws.Range("C9").Value = "31.01."
$ws.Range("D9").Value = "MITGLIEDSBEITRAG ZEUS BODYPOWER"
$ws.Range("E9").Value = "25,27-"

# Row 10
$ws.Range("B10").Value = "02.02."
$ws.Range("C10").Value = "03.02."
$ws.Range("D10").Value = "AMAZON.DE MKTPLC EU AWHPYH"
$ws.Range("E10").Value = "204,18-"

# Row 11: previously blank, now a new transaction row.
# Copy the style of row 10 (B:D) into row 11, and match E-column style to the
# other amount cells (style used by E6:E10) rather than the old blank style.
$ws.Range("B10:D10").Copy()
$ws.Range("B11:D11").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("E10").Copy()
$ws.Range("E11").PasteSpecial(-4122) # xlPasteFormats

$ws.Range("B11").Value = "03.02."
$ws.Range("C11").Value = "04.02."
$ws.Range("D11").Value = "BEITRAG Allianz SE K-18086548"
$ws.Range("E11").Value = "55,75-"

# Closing balance line
$ws.Range("D12").Value = "KONTOSTAND AM 07.02.2024"
$ws.Range("E12").Value = "458,51-"

# Next statement date
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 16.02.2024"
